$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 834, shifting existing rows 834:875 down to 835:876
$ws.Rows(834).Insert()

# Populate the newly inserted row 834 with the new record.
# Force column A to be stored as text (not auto-parsed into a date serial)
# by temporarily formatting as Text, then clear the formatting afterwards so
# the cell ends up with the same (default) style as its neighbours.
$ws.Range("A834").NumberFormat = "@"
$ws.Range("A834").Value = "2026/02/22"
$ws.Range("A834").ClearFormats()

$ws.Range("B834").Value = "日"
$ws.Range("C834").Value = 16
$ws.Range("D834").Value = 25
